$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'95.338.65"
$ws.Range("E2").Value2 = "  +0.29%  "

$ws.Range("D3").Value2 = "'3.561.42"
$ws.Range("E3").Value2 = "  +0.45%  "

$ws.Range("E4").Value2 = "  +0.23%  "

$ws.Range("D5").Value2 = "'235.16"
$ws.Range("E5").Value2 = "  -1.82%  "

$ws.Range("D6").Value2 = "'649.17"
$ws.Range("E6").Value2 = "  +2.09%  "

$ws.Range("E7").Value2 = "  -0.37%  "

$ws.Range("D8").Value2 = "'0.396"
$ws.Range("E8").Value2 = "  -1.18%  "

$ws.Range("E9").Value2 = "  +0.16%  "

$ws.Range("D10").Value2 = "'0.978"
$ws.Range("E10").Value2 = "  -2.88%  "

$ws.Range("D11").Value2 = "'3.559.67"
$ws.Range("E11").Value2 = "  +0.64%  "

$ws.Range("E12").Value2 = "  -0.22%  "

$ws.Range("D13").Value2 = "'41.99"
$ws.Range("E13").Value2 = "  -4.42%  "

$ws.Range("D14").Value2 = "'6.50"
$ws.Range("E14").Value2 = "  +3.15%  "

$ws.Range("D15").Value2 = "'4.241.02"
$ws.Range("E15").Value2 = "  +0.66%  "

$ws.Range("D16").Value2 = "'95.355.65"
$ws.Range("E16").Value2 = "  +0.69%  "

$ws.Range("E17").Value2 = "  -0.95%  "

$ws.Range("D18").Value2 = "'3.555.60"
$ws.Range("E18").Value2 = "  +0.26%  "

$ws.Range("E19").Value2 = "  -5.71%  "

$ws.Range("D20").Value2 = "'12.65"
$ws.Range("E20").Value2 = "  -2.50%  "

$ws.Range("D21").Value2 = "'17.63"
$ws.Range("E21").Value2 = "  -2.45%  "

$ws.Range("D22").Value2 = "'3.46"
$ws.Range("E22").Value2 = "  +0.87%  "

$ws.Range("D23").Value2 = "'503.38"
$ws.Range("E23").Value2 = "  -2.62%  "

$ws.Range("D24").Value2 = "'0.470"
$ws.Range("E24").Value2 = "  -5.33%  "

$ws.Range("E25").Value2 = "  -0.94%  "

$ws.Range("E26").Value2 = "  -2.98%  "

$ws.Range("D27").Value2 = "'91.37"
$ws.Range("E27").Value2 = "  -4.94%  "

$ws.Range("D28").Value2 = "'12.38"
$ws.Range("E28").Value2 = "  +0.38%  "

$ws.Range("D29").Value2 = "'3.753.45"
$ws.Range("E29").Value2 = "  +0.86%  "

$ws.Range("D30").Value2 = "'3.03"
$ws.Range("E30").Value2 = "  -0.71%  "

$ws.Range("D31").Value2 = "'0.999"
$ws.Range("E31").Value2 = "  -0.11%  "

$ws.Range("D32").Value2 = "'11.18"
$ws.Range("E32").Value2 = "  -3.68%  "

$ws.Range("E34").Value2 = "  +0.54%  "

$ws.Range("E35").Value2 = "  -2.76%  "

$ws.Range("D36").Value2 = "'31.77"
$ws.Range("E36").Value2 = "  +5.05%  "

$ws.Range("D37").Value2 = "'0.555"
$ws.Range("E37").Value2 = "  -1.95%  "

$ws.Range("D38").Value2 = "'8.17"
$ws.Range("E38").Value2 = "  +7.11%  "

$ws.Range("D39").Value2 = "'559.53"
$ws.Range("E39").Value2 = "  -5.05%  "

$ws.Range("D40").Value2 = "'1.51"
$ws.Range("E40").Value2 = "  +3.80%  "

$ws.Range("E41").Value2 = "  -0.02%  "

$ws.Range("D42").Value2 = "'0.149"
$ws.Range("E42").Value2 = "  -1.18%  "

$ws.Range("D43").Value2 = "'0.896"
$ws.Range("E43").Value2 = "  -3.91%  "

$ws.Range("E44").Value2 = "  +2.01%  "

$ws.Range("D45").Value2 = "'34.74"
$ws.Range("E45").Value2 = "  +34.40%  "

$ws.Range("D46").Value2 = "'2.28"
$ws.Range("E46").Value2 = "  +4.82%  "

$ws.Range("E47").Value2 = "  -1.01%  "

$ws.Range("D48").Value2 = "'5.59"
$ws.Range("E48").Value2 = "  -0.13%  "

$ws.Range("D49").Value2 = "'0.0407"
$ws.Range("E49").Value2 = "  -4.57%  "

$ws.Range("D50").Value2 = "'3.57"
$ws.Range("E50").Value2 = "  +0.14%  "

$ws.Range("D51").Value2 = "'53.33"
$ws.Range("E51").Value2 = "  -1.04%  "
